$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (old rows 8-10, the MuSCs-as-sender block)
$ws.Range("A8:T10").EntireRow.Delete()

# Refresh the remaining six data rows (now rows 2-7) with the updated TPM-based
# NATMI values. Column order: A Sending cluster, B Ligand symbol, C Receptor
# symbol, D Target cluster, E..T the numeric NATMI statistics.

$ws.Cells.Item(2,1).Value  = "ECs"
$ws.Cells.Item(2,2).Value  = "Ntf3"
$ws.Cells.Item(2,3).Value  = "Ntrk3"
$ws.Cells.Item(2,4).Value  = "FAPs"
$ws.Cells.Item(2,5).Value  = 3
$ws.Cells.Item(2,6).Value  = 1
$ws.Cells.Item(2,7).Value  = 9.736532333333333
$ws.Cells.Item(2,8).Value  = 29.209597
$ws.Cells.Item(2,9).Value  = 0.3545698647072128
$ws.Cells.Item(2,10).Value = 0.3545698647072129
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.468673666666667
$ws.Cells.Item(2,14).Value = 7.406021
$ws.Cells.Item(2,15).Value = 0.635345274347677
$ws.Cells.Item(2,16).Value = 0.635345274347677
$ws.Cells.Item(2,17).Value = 24.03632097594856
$ws.Cells.Item(2,18).Value = 216.326888783537
$ws.Cells.Item(2,19).Value = 0.2252742879678228
$ws.Cells.Item(2,20).Value = 0.2252742879678229

$ws.Cells.Item(3,1).Value  = "ECs"
$ws.Cells.Item(3,2).Value  = "Ntf3"
$ws.Cells.Item(3,3).Value  = "Ntrk3"
$ws.Cells.Item(3,4).Value  = "MuSCs"
$ws.Cells.Item(3,5).Value  = 3
$ws.Cells.Item(3,6).Value  = 1
$ws.Cells.Item(3,7).Value  = 9.736532333333333
$ws.Cells.Item(3,8).Value  = 29.209597
$ws.Cells.Item(3,9).Value  = 0.3545698647072128
$ws.Cells.Item(3,10).Value = 0.3545698647072129
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.416888666666667
$ws.Cells.Item(3,14).Value = 4.250666
$ws.Cells.Item(3,15).Value = 0.364654725652323
$ws.Cells.Item(3,16).Value = 0.364654725652323
$ws.Cells.Item(3,17).Value = 13.79558231573355
$ws.Cells.Item(3,18).Value = 124.160240841602
$ws.Cells.Item(3,19).Value = 0.12929557673939
$ws.Cells.Item(3,20).Value = 0.12929557673939

$ws.Cells.Item(4,1).Value  = "FAPs"
$ws.Cells.Item(4,2).Value  = "Ntf3"
$ws.Cells.Item(4,3).Value  = "Ntrk3"
$ws.Cells.Item(4,4).Value  = "FAPs"
$ws.Cells.Item(4,5).Value  = 3
$ws.Cells.Item(4,6).Value  = 1
$ws.Cells.Item(4,7).Value  = 10.17625966666667
$ws.Cells.Item(4,8).Value  = 30.528779
$ws.Cells.Item(4,9).Value  = 0.37058316962423
$ws.Cells.Item(4,10).Value = 0.37058316962423
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.468673666666667
$ws.Cells.Item(4,14).Value = 7.406021
$ws.Cells.Item(4,15).Value = 0.635345274347677
$ws.Cells.Item(4,16).Value = 0.635345274347677
$ws.Cells.Item(4,17).Value = 25.12186426426211
$ws.Cells.Item(4,18).Value = 226.096778378359
$ws.Cells.Item(4,19).Value = 0.2354482655735381
$ws.Cells.Item(4,20).Value = 0.2354482655735382

$ws.Cells.Item(5,1).Value  = "FAPs"
$ws.Cells.Item(5,2).Value  = "Ntf3"
$ws.Cells.Item(5,3).Value  = "Ntrk3"
$ws.Cells.Item(5,4).Value  = "MuSCs"
$ws.Cells.Item(5,5).Value  = 3
$ws.Cells.Item(5,6).Value  = 1
$ws.Cells.Item(5,7).Value  = 10.17625966666667
$ws.Cells.Item(5,8).Value  = 30.528779
$ws.Cells.Item(5,9).Value  = 0.37058316962423
$ws.Cells.Item(5,10).Value = 0.37058316962423
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.416888666666667
$ws.Cells.Item(5,14).Value = 4.250666
$ws.Cells.Item(5,15).Value = 0.364654725652323
$ws.Cells.Item(5,16).Value = 0.364654725652323
$ws.Cells.Item(5,17).Value = 14.41862699075711
$ws.Cells.Item(5,18).Value = 129.767642916814
$ws.Cells.Item(5,19).Value = 0.1351349040506919
$ws.Cells.Item(5,20).Value = 0.1351349040506919

$ws.Cells.Item(6,1).Value  = "MuSCs"
$ws.Cells.Item(6,2).Value  = "Ntf3"
$ws.Cells.Item(6,3).Value  = "Ntrk3"
$ws.Cells.Item(6,4).Value  = "FAPs"
$ws.Cells.Item(6,5).Value  = 3
$ws.Cells.Item(6,6).Value  = 1
$ws.Cells.Item(6,7).Value  = 7.547331666666667
$ws.Cells.Item(6,8).Value  = 22.641995
$ws.Cells.Item(6,9).Value  = 0.2748469656685572
$ws.Cells.Item(6,10).Value = 0.2748469656685572
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.468673666666667
$ws.Cells.Item(6,14).Value = 7.406021
$ws.Cells.Item(6,15).Value = 0.635345274347677
$ws.Cells.Item(6,16).Value = 0.635345274347677
$ws.Cells.Item(6,17).Value = 18.63189893909945
$ws.Cells.Item(6,18).Value = 167.687090451895
$ws.Cells.Item(6,19).Value = 0.174622720806316
$ws.Cells.Item(6,20).Value = 0.174622720806316

$ws.Cells.Item(7,1).Value  = "MuSCs"
$ws.Cells.Item(7,2).Value  = "Ntf3"
$ws.Cells.Item(7,3).Value  = "Ntrk3"
$ws.Cells.Item(7,4).Value  = "MuSCs"
$ws.Cells.Item(7,5).Value  = 3
$ws.Cells.Item(7,6).Value  = 1
$ws.Cells.Item(7,7).Value  = 7.547331666666667
$ws.Cells.Item(7,8).Value  = 22.641995
$ws.Cells.Item(7,9).Value  = 0.2748469656685572
$ws.Cells.Item(7,10).Value = 0.2748469656685572
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.416888666666667
$ws.Cells.Item(7,14).Value = 4.250666
$ws.Cells.Item(7,15).Value = 0.364654725652323
$ws.Cells.Item(7,16).Value = 0.364654725652323
$ws.Cells.Item(7,17).Value = 10.69372870207444
$ws.Cells.Item(7,18).Value = 96.24355831867
$ws.Cells.Item(7,19).Value = 0.1002242448622411
$ws.Cells.Item(7,20).Value = 0.1002242448622411
